# NATMI "Hbegf-Erbb2" LR-pairs sheet refresh ("update scripts wuth new tpm").
#
# The ligand/receptor summary table (rows 2..N) is the full cross-join of
# Sending cluster (col A) x Target cluster (col D) over the 4 clusters
# {ECs, FAPs, MuSCs, Resolving-Mac}, with the ligand always Hbegf (col B)
# and the receptor always Erbb2 (col C). The previous run only emitted 3
# of the 4 target clusters per sender (12 rows, missing "Resolving-Mac" as
# a target); the refreshed TPM run emits the complete 4x4 grid (16 rows)
# with updated expression/specificity numbers throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Each row: Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average/total
# expression value, Ligand derived specificity (avg/total), Receptor-expressing
# cells, Receptor detection rate, Receptor average/total expression value,
# Receptor derived specificity (avg/total), Edge average/total expression
# weight, Edge average/total expression derived specificity.
$data = @(
    ,@("ECs", "Hbegf", "Erbb2", "ECs", 3, 1, 8.082188, 24.246564, 0.2755301789948819, 0.2755301789948819, 3, 1, 3.020000333333333, 9.060001, 0.291481777372034, 0.291481777372034, 24.40821045406267, 219.673894086564, 0.08031202629306283, 0.08031202629306285)
    ,@("ECs", "Hbegf", "Erbb2", "FAPs", 3, 1, 8.082188, 24.246564, 0.2755301789948819, 0.2755301789948819, 3, 1, 4.076388666666666, 12.229166, 0.3934413518781783, 0.3934413518781784, 32.94613956506933, 296.515256085624, 0.1084049661069828, 0.1084049661069828)
    ,@("ECs", "Hbegf", "Erbb2", "MuSCs", 3, 1, 8.082188, 24.246564, 0.2755301789948819, 0.2755301789948819, 3, 1, 3.229698, 9.689094000000001, 0.311721195201271, 0.3117211952012711, 26.103026419224, 234.9272377730161, 0.0858885967103047, 0.08588859671030474)
    ,@("ECs", "Hbegf", "Erbb2", "Resolving-Mac", 3, 1, 8.082188, 24.246564, 0.2755301789948819, 0.2755301789948819, 1, 0.3333333333333333, 0.03476766666666667, 0.104303, 0.003355675548516525, 0.003355675548516525, 0.2809988183213334, 2.528989364892, 0.0009245898845315065, 0.0009245898845315068)
    ,@("FAPs", "Hbegf", "Erbb2", "ECs", 3, 1, 13.29805733333333, 39.894172, 0.4533445791334642, 0.4533445791334642, 3, 1, 3.020000333333333, 9.060001, 0.291481777372034, 0.291481777372034, 40.16013757935244, 361.441238214172, 0.1321416836877989, 0.1321416836877989)
    ,@("FAPs", "Hbegf", "Erbb2", "FAPs", 3, 1, 13.29805733333333, 39.894172, 0.4533445791334642, 0.4533445791334642, 3, 1, 4.076388666666666, 12.229166, 0.3934413518781783, 0.3934413518781784, 54.20805020228355, 487.8724518205519, 0.1783645040809139, 0.1783645040809139)
    ,@("FAPs", "Hbegf", "Erbb2", "MuSCs", 3, 1, 13.29805733333333, 39.894172, 0.4533445791334642, 0.4533445791334642, 3, 1, 3.229698, 9.689094000000001, 0.311721195201271, 0.3117211952012711, 42.948709173352, 386.538382560168, 0.1413171140455006, 0.1413171140455007)
    ,@("FAPs", "Hbegf", "Erbb2", "Resolving-Mac", 3, 1, 13.29805733333333, 39.894172, 0.4533445791334642, 0.4533445791334642, 1, 0.3333333333333333, 0.03476766666666667, 0.104303, 0.003355675548516525, 0.003355675548516525, 0.4623424246795556, 4.161081822116, 0.00152127731925068, 0.001521277319250681)
    ,@("MuSCs", "Hbegf", "Erbb2", "ECs", 3, 1, 5.789497666666667, 17.368493, 0.1973699855023315, 0.1973699855023315, 3, 1, 3.020000333333333, 9.060001, 0.291481777372034, 0.291481777372034, 17.48428488316589, 157.358563948493, 0.05752975417411216, 0.05752975417411217)
    ,@("MuSCs", "Hbegf", "Erbb2", "FAPs", 3, 1, 5.789497666666667, 17.368493, 0.1973699855023315, 0.1973699855023315, 3, 1, 4.076388666666666, 12.229166, 0.3934413518781783, 0.3934413518781784, 23.60024267409311, 212.402184066838, 0.07765351391621375, 0.07765351391621376)
    ,@("MuSCs", "Hbegf", "Erbb2", "MuSCs", 3, 1, 5.789497666666667, 17.368493, 0.1973699855023315, 0.1973699855023315, 3, 1, 3.229698, 9.689094000000001, 0.311721195201271, 0.3117211952012711, 18.698329035038, 168.284961315342, 0.0615244077776443, 0.06152440777764432)
    ,@("MuSCs", "Hbegf", "Erbb2", "Resolving-Mac", 3, 1, 5.789497666666667, 17.368493, 0.1973699855023315, 0.1973699855023315, 1, 0.3333333333333333, 0.03476766666666667, 0.104303, 0.003355675548516525, 0.003355675548516525, 0.2012873250421111, 1.811585925379, 0.0006623096343612348, 0.0006623096343612349)
    ,@("Resolving-Mac", "Hbegf", "Erbb2", "ECs", 3, 1, 2.163479333333334, 6.490438, 0.0737552563693224, 0.0737552563693224, 3, 1, 3.020000333333333, 9.060001, 0.291481777372034, 0.291481777372034, 6.533708307826445, 58.803374770438, 0.02149831321706012, 0.02149831321706012)
    ,@("Resolving-Mac", "Hbegf", "Erbb2", "FAPs", 3, 1, 2.163479333333334, 6.490438, 0.0737552563693224, 0.0737552563693224, 3, 1, 4.076388666666666, 12.229166, 0.3934413518781783, 0.3934413518781784, 8.819182634967556, 79.37264371470799, 0.02901836777406782, 0.02901836777406783)
    ,@("Resolving-Mac", "Hbegf", "Erbb2", "MuSCs", 3, 1, 2.163479333333334, 6.490438, 0.0737552563693224, 0.0737552563693224, 3, 1, 3.229698, 9.689094000000001, 0.311721195201271, 0.3117211952012711, 6.987384875908002, 62.88646388317201, 0.02299107666782133, 0.02299107666782134)
    ,@("Resolving-Mac", "Hbegf", "Erbb2", "Resolving-Mac", 3, 1, 2.163479333333334, 6.490438, 0.0737552563693224, 0.0737552563693224, 1, 0.3333333333333333, 0.03476766666666667, 0.104303, 0.003355675548516525, 0.003355675548516525, 0.07521912830155557, 0.6769721547140001, 0.0002474987103731029, 0.0002474987103731029)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$r").Value = $rowValues[$j]
    }
}
